$d = $word.ActiveDocument

# Map of exact paragraph text (English) -> replacement text (Italian).
# Using exact full-paragraph matches keeps this safe from accidentally
# touching the "English" inside the language-switcher hyperlink at the
# top of the document (that paragraph's full text is the longer
# "English / Portuguese / French / Thai / Vietnamese / Spanish" line).
$replacements = @{
    "English" = "Inglese";
    "Don’t forget to send your documents" = "Non dimenticarti di inviare i documenti";
    "If you have any questions, please contact your country manager." = "Per qualsiasi domanda, contatta il tuo country manager.";
    "We look forward to seeing you there!" = "Non vediamo l'ora di incontrarti!"
}

foreach ($p in $d.Paragraphs) {
    $rng = $p.Range
    $current = $rng.Text
    # Paragraph.Range.Text includes the trailing paragraph mark; trim it
    # off before comparing against the target phrase.
    $currentTrimmed = $current.TrimEnd([char]13, [char]7)

    if ($replacements.ContainsKey($currentTrimmed)) {
        $newText = $replacements[$currentTrimmed]
        $found = $rng.Find.Execute($currentTrimmed, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Text = $newText
        }
    }
}
